# Naive Bayes armazenando valores
# Adds 24 new "Naive Bayes" result rows (20-43) to Sheet1 and fixes the
# "Atividade" column for the existing "PCA ..." SVC rows (14-19) so they
# hold the number 7 instead of the text "7".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A14:A19: was stored as text "7", should be the number 7 --------
$ws.Range("A14:A19").Value = 7

# --- New data rows 20-43: Naive Bayes results ----------------------------
$rows = @(
    @(20, 6,   "Naive Bayes", "Branch and Bound Desbalanceado",  0.6, "Accuracy",  0.9312448700410396),
    @(21, 6,   "Naive Bayes", "Branch and Bound Desbalanceado",  0.6, "Recall",    0),
    @(22, 6,   "Naive Bayes", "Branch and Bound Desbalanceado",  0.6, "Precision", 0),
    @(23, 6,   "Naive Bayes", "Branch and Bound Balanceado",     0.3, "Accuracy",  0.6901094391244869),
    @(24, 6,   "Naive Bayes", "Branch and Bound Balanceado",     0.3, "Recall",    0.2166666666666667),
    @(25, 6,   "Naive Bayes", "Branch and Bound Balanceado",     0.3, "Precision", 0.04927119294207902),
    @(26, 6,   "Naive Bayes", "Dataset Completo Desbalanceado",  0.6, "Accuracy",  0.9556771545827634),
    @(27, 6,   "Naive Bayes", "Dataset Completo Desbalanceado",  0.6, "Recall",    0.7166666666666666),
    @(28, 6,   "Naive Bayes", "Dataset Completo Desbalanceado",  0.6, "Precision", 0.6753968253968254),
    @(29, 6,   "Naive Bayes", "Dataset Completo Balanceado",     0.8, "Accuracy",  0.9393296853625172),
    @(30, 6,   "Naive Bayes", "Dataset Completo Balanceado",     0.8, "Recall",    0.8766666666666666),
    @(31, 6,   "Naive Bayes", "Dataset Completo Balanceado",     0.8, "Precision", 0.547142857142857),
    @(32, 6,   "Naive Bayes", "PCA Desbalanceado",                1,  "Accuracy",  0.9335567715458277),
    @(33, 6,   "Naive Bayes", "PCA Desbalanceado",                1,  "Recall",    0.11),
    @(34, 6,   "Naive Bayes", "PCA Desbalanceado",                1,  "Precision", 0.2166666666666666),
    @(35, 6,   "Naive Bayes", "PCA Balanceado",                   1,  "Accuracy",  0.9416963064295485),
    @(36, 6,   "Naive Bayes", "PCA Balanceado",                   1,  "Recall",    0.64),
    @(37, 6,   "Naive Bayes", "PCA Balanceado",                   1,  "Precision", 0.5704761904761905),
    @(38, "6", "Naive Bayes", "ReliefF Desbalanceado",           0.4, "Accuracy",  0.9627222982216141),
    @(39, "6", "Naive Bayes", "ReliefF Desbalanceado",           0.4, "Recall",    0.6966666666666667),
    @(40, "6", "Naive Bayes", "ReliefF Desbalanceado",           0.4, "Precision", 0.7538095238095239),
    @(41, "6", "Naive Bayes", "ReliefF Balanceado",              0.1, "Accuracy",  0.9568946648426813),
    @(42, "6", "Naive Bayes", "ReliefF Balanceado",              0.1, "Recall",    0.8766666666666666),
    @(43, "6", "Naive Bayes", "ReliefF Balanceado",              0.1, "Precision", 0.6523015873015873)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]

    $aCell = $ws.Cells.Item($r, 1)
    if ($row[1] -is [string]) {
        # Store a numeric-looking value as literal text (no style change):
        # write it as a formula returning the text, then collapse the
        # formula down to its cached value via copy / paste-values.
        $aCell.Formula = '="' + $row[1] + '"'
        $aCell.Copy()
        $aCell.PasteSpecial(-4163) # xlPasteValues
        $excel.CutCopyMode = $false
    } else {
        $aCell.Value = $row[1]
    }
}

Write-Output "ok"
